$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-5: 45183 -> 45184
foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
